$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly measured data for columns D (dry weight fraction) and E, rows 42-69,
# filling in the gaps so the full sample set (rows 2-81) is consistent.
$data = @{
    42 = @(0.477,  0.3031)
    43 = @(0.5437, 0.3444)
    44 = @(0.4877, 0.3094)
    45 = @(0.5822, 0.3673)
    46 = @(0.5148, 0.3313)
    47 = @(0.5782, 0.3711)
    48 = @(0.5753, 0.369)
    49 = @(0.5604, 0.3488)
    50 = @(0.5893, 0.3735)
    51 = @(0.3472, 0.2161)
    52 = @(0.567,  0.367)
    53 = @(0.5194, 0.3298)
    54 = @(0.5326, 0.3372)
    55 = @(0.4568, 0.2881)
    56 = @(0.4991, 0.3079)
    57 = @(0.571,  0.3564)
    58 = @(0.5548, 0.3537)
    59 = @(0.5764, 0.3638)
    60 = @(0.5826, 0.3675)
    61 = @(0.521,  0.3264)
    62 = @(0.535,  0.3299)
    63 = @(0.5683, 0.3523)
    64 = @(0.5341, 0.335)
    65 = @(0.5458, 0.3416)
    66 = @(0.5227, 0.3285)
    67 = @(0.528,  0.3313)
    68 = @(0.5045, 0.3138)
    69 = @(0.474,  0.2975)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
}

# Move the view / active selection to where the newly-filled rows are.
$ws.Range("E70").Select()
